# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the per-locale report sheets (zh-cn, de-de) to reflect the
# freshly generated handback report.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: Correspond Handoff Datetime (col D) / Correspond Handback DateTime (col G)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-25 06:24:10"
$wsZhCn.Range("D3").Value = "2016-02-25 06:24:10"
$wsZhCn.Range("G2").Value = "2016-02-25 06:24:59"
$wsZhCn.Range("G3").Value = "2016-02-25 06:24:59"

# de-de sheet: Correspond Handoff Datetime (col D) / Correspond Handback DateTime (col G)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-25 06:24:22"
$wsDeDe.Range("D3").Value = "2016-02-25 06:24:22"
$wsDeDe.Range("G2").Value = "2016-02-25 06:25:22"
$wsDeDe.Range("G3").Value = "2016-02-25 06:25:22"
